$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.95"
$ws.Range("E2").Value = "'0.21%"
$ws.Range("D3").Value = "'26.92"
$ws.Range("E3").Value = "'-0.21%"
$ws.Range("D4").Value = "'4.675"
$ws.Range("E4").Value = "'-1.02%"
$ws.Range("D5").Value = "'0.05891"
$ws.Range("E5").Value = "'-0.77%"
$ws.Range("D6").Value = "'6.610"
$ws.Range("D7").Value = "'0.8499"
$ws.Range("E7").Value = "'-2.09%"
$ws.Range("D8").Value = "'0.9206"
$ws.Range("E8").Value = "'-3.47%"
$ws.Range("D9").Value = "'0.1375"
$ws.Range("E9").Value = "'-2.07%"
$ws.Range("D10").Value = "'0.04226"
$ws.Range("E10").Value = "'6.84%"
$ws.Range("D11").Value = "'0.06993"
$ws.Range("E11").Value = "'-2.44%"
$ws.Range("D12").Value = "'0.03052"
$ws.Range("E12").Value = "'-4.21%"
$ws.Range("D13").Value = "'0.09109"
$ws.Range("E13").Value = "'-1.57%"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'-0.92%"
$ws.Range("B15").Value = "'One"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006040"
$ws.Range("E15").Value = "'-0.50%"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006019"
$ws.Range("E16").Value = "'-0.93%"
$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.469"
$ws.Range("E17").Value = "'-0.40%"
$ws.Range("B18").Value = "'GateToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.166"
$ws.Range("E18").Value = "'-1.06%"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.177"
$ws.Range("E19").Value = "'-1.89%"
$ws.Range("D20").Value = "'0.3046"
$ws.Range("E20").Value = "'-2.82%"
$ws.Range("E21").Value = "'-0.23%"
$ws.Range("D22").Value = "'3.919"
$ws.Range("E22").Value = "'3.04%"
$ws.Range("D23").Value = "'0.04235"
$ws.Range("E23").Value = "'0.35%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.07%"
$ws.Range("D25").Value = "'0.004293"
$ws.Range("E25").Value = "'-4.51%"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("E27").Value = "'-21.36%"
$ws.Range("D40").Value = "'0.03781"
$ws.Range("E40").Value = "'-1.12%"
$ws.Range("D41").Value = "'0.006239"
$ws.Range("E41").Value = "'7.40%"
$ws.Range("D42").Value = "'0.1099"
$ws.Range("E42").Value = "'-0.01%"
$ws.Range("D43").Value = "'0.002430"
$ws.Range("E43").Value = "'15.42%"
$ws.Range("D44").Value = "'0.01416"
$ws.Range("E44").Value = "'34.02%"
$ws.Range("D45").Value = "'0.00005358"
$ws.Range("E45").Value = "'-2.59%"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E47").Value = "'-48.60%"
$ws.Range("E48").Value = "'10,463.88%"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E50").Value = "'0.01%"
